$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "extra" rows that don't belong in the original test data are
# old row 20 (TAC00066440 / Acetate#UG-8940 ...) and old row 4
# (TAC11187850 / Acetate#UN-0041 ...). Delete the higher-numbered row
# first so the lower row index is still valid afterwards.
$ws.Rows(20).Delete()
$ws.Rows(4).Delete()

# Re-apply the AutoFilter over the new (smaller) used range and restore
# the value-filter on column D (field 4) to TAC00066440.
$ws.AutoFilterMode = $false
$ws.Range("A1:Q18").AutoFilter(4, @("TAC00066440"), 7)

# Re-apply the descending sort on column D over the whole filtered range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D1:D18"), 0, 2, 0, 0)
$ws.Sort.SetRange($ws.Range("A1:Q18"))
$ws.Sort.Apply()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new AutoFilter range.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$Q`$18"

# Move the selection down to row 19 (first row below the new data), like
# the original edit's saved selection/scroll state.
$ws.Rows(19).Select()
